# Mise a jour de l'application
# Adds 10 new wellness entries (rows 509-518, all dated 2025-10-28 / serial 45958)
# below the existing data table, reusing the formatting of the last data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend the table formatting down by copying the last populated row's
#    look (date format in A, data font in B:H) into the 10 new rows.
$ws.Rows(508).Copy()
$ws.Rows("509:518").Insert()

# 2) New values for the appended rows.
$newRows = @(
    @{ Row=509; Player="Kamal Bafounta";  C=70; D=5; E=3; F=2; G="Genou";        H=5 },
    @{ Row=510; Player="Levy Ndoutoume";  C=70; D=7; E=7; F=2; G="Ischio";       H=5 },
    @{ Row=511; Player="Romain Thunet";   C=70; D=6; E=6; F=2; G="Synthétique "; H=9 },
    @{ Row=512; Player="Omar Benyounes";  C=70; D=6; E=7; F=2; G="Coup pied";    H=1 },
    @{ Row=513; Player="Yoann Martelat";  C=70; D=5; E=5; F=5; G="Genou";        H=6 },
    @{ Row=514; Player="Malik Boussaid";  C=70; D=2; E=0; F=0; G=$null;          H=10 },
    @{ Row=515; Player="Karim Belmahi";   C=70; D=6; E=6; F=0; G=$null;          H=10 },
    @{ Row=516; Player="Mattheo Haon";    C=70; D=7; E=6; F=0; G=$null;          H=10 },
    @{ Row=517; Player="Sofiane Belle";   C=70; D=6; E=4; F=0; G=$null;          H=8 },
    @{ Row=518; Player="Karahali Souaré"; C=70; D=4; E=6; F=7; G="Cheville";     H=2 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 45958
    $ws.Cells.Item($row, 2).Value = $r.Player
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    if ($r.G) {
        $ws.Cells.Item($row, 7).Value = $r.G
    }
    $ws.Cells.Item($row, 8).Value = $r.H
}

# 3) Column I: Charge = Volume * Intensite, filled down as shared formulas
#    (matches the original I452:I508 shared formula pattern, continued then
#    split into a second shared group for the last three new rows).
$ws.Range("I509:I515").Formula = "=C509*D509"
$ws.Range("I516:I518").Formula = "=C516*D516"

# 4) Update the visible window / active selection to reflect the new bottom
#    of the sheet, same as what Excel records after scrolling to the new rows.
$ws.Range("A486").Select()
$excel.ActiveWindow.ScrollRow = 486
$ws.Range("K513").Select()
